$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 01:46:58"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6 updates
$ws1.Range("A6").Value = "01:46:58"
$ws1.Range("B6").Value = "01:57"
$ws1.Range("D6").Value = 11

# Row 7 updates
$ws1.Range("A7").Value = "01:46:58"
$ws1.Range("D7").Value = 71

# New row 8
$ws1.Range("A8").Value = "01:46:58"
$ws1.Range("B8").Value = "03:04"
$ws1.Range("C8").Value = "15_ABASTO"
$ws1.Range("D8").Value = 78
$ws1.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:46:58"
$ws2.Range("A6").Value = "01:46:58"
$ws2.Range("D6").Value = 71

# ---------------------------------------------------------------
# Sheet 3: "6203-6173"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:46:58"
